$d = $word.ActiveDocument

# Locate the paragraph that hosts the "promptString" field (field code
# runs + fldChar begin/end) and rewrite it so the field is represented
# as literal text runs instead (TokenIteratorFieldRewriterSplit style).
$targetPara = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text -match "promptString") {
        $targetPara = $f.Code.Paragraphs(1)
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate the promptString field paragraph"
}

$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve">''A </w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>String</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve"> please''.</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>p</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>rompt</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>String</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>()</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p>'

$targetPara.Range.InsertXML($fragment)
